# "iPad synth version finished"
# - Battery LiPo capacity corrected (660mA -> 6600mA) and quantity doubled
# - Several order quantities doubled (Charger, SHIM, Power supply, IMU breakout, Ethernet cable)
# - D4 / D5 hyperlink text & targets swapped
# - Selection moved to B33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity changes (column A) ---
$ws.Range("A4").Value  = 2    # Battery LiPo 6600mA
$ws.Range("A5").Value  = 2    # Charger - PowerBoost 1000 Charger
$ws.Range("A7").Value  = 4    # Pimoroni ONOff SHIM for Raspberry Pi
$ws.Range("A12").Value = 2    # Power supply
$ws.Range("A27").Value = 2    # LDM9D51 IMU breakout
$ws.Range("A28").Value = 2    # Ethernet panel mount cable

# --- Battery description text (row 4) ---
$ws.Range("B4").Value = "Battery LiPo 6600mA"

# --- Rebuild the hyperlinks: delete everything then re-add in the updated order,
#     swapping D4 / D5's displayed text + targets ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D12"), "https://www.adafruit.com/product/1994")
$ws.Hyperlinks.Add($ws.Range("D9"),  "https://www.adafruit.com/product/1994")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.adafruit.com/product/1445")
$ws.Hyperlinks.Add($ws.Range("D8"),  "https://www.adafruit.com/product/559")
$ws.Hyperlinks.Add($ws.Range("D6"),  "https://www.adafruit.com/product/2225")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.adafruit.com/product/2046")
$ws.Hyperlinks.Add($ws.Range("C23"), "https://www.mouser.ca/ProductDetail/485-1988")
$ws.Hyperlinks.Add($ws.Range("C24"), "https://www.mouser.ca/ProductDetail/485-937")
$ws.Hyperlinks.Add($ws.Range("C26"), "https://www.mouser.ca/ProductDetail/Qualtek/3021007-06?qs=sGAEpiMZZMsgIz308WEU047hkcLyYWNnYAdW6L8LtYQ%3d")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://www.mouser.ca/ProductDetail/858-P160KNP0QC20A10K")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://www.mouser.ca/ProductDetail/485-1988")
$ws.Hyperlinks.Add($ws.Range("C22"), "https://www.mouser.ca/ProductDetail/81-7BB-20-6L0")
$ws.Hyperlinks.Add($ws.Range("C25"), "https://www.mouser.ca/ProductDetail/485-1852")
$ws.Hyperlinks.Add($ws.Range("C27"), "https://www.mouser.ca/ProductDetail/485-3387")
$ws.Hyperlinks.Add($ws.Range("C28"), "https://www.mouser.ca/ProductDetail/485-909")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://www.mouser.ca/ProjectManager/ProjectDetail.aspx?AccessID=2ecda1393f")
$ws.Hyperlinks.Add($ws.Range("D7"),  "https://www.adafruit.com/product/3581")

# D4 / D5: new hyperlink objects, text & target swapped relative to original
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.adafruit.com/product/353", "", "", "https://www.adafruit.com/product/353")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.adafruit.com/product/2465", "", "", "https://www.adafruit.com/product/2465")

# --- Move the active selection to B33 ---
$ws.Range("B33").Select()
